$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    32.99293476771223,
    39.81795260052328,
    32.8466741743693,
    30.89634255962685,
    30.76680808256089,
    26.39873364736599,
    24.86969599478817,
    26.38991189248756,
    23.94022531284983,
    27.25960600151734,
    24.75452262489402,
    20.46655330024399,
    21.05394258799927,
    22.20226048738496,
    19.46990463595043,
    17.26670102626844,
    21.60091681413137,
    15.00795634773498,
    20.09772569057796,
    19.87328033220207,
    17.24699887520443,
    18.39964834137456,
    16.9082392868022,
    16.73228708559159,
    13.76641543155043,
    17.51973174560494,
    18.69396516350307,
    14.47722281418032,
    13.18337501373389,
    10.57188800643623,
    15.11393900861776,
    18.11137515097056,
    17.68118765067433,
    17.31761244851236,
    14.39128494512975,
    13.73328375738535,
    14.79663574215394,
    15.40097082537302,
    12.75473845065363,
    12.72654919580208,
    15.51420423065211,
    15.17978664506111,
    15.40841250596976,
    15.90846253914941,
    16.43305033557174,
    15.6487362404888,
    13.85595256237559,
    12.96146455252259,
    15.12808471808341,
    11.78453101333216,
    13.48623860161457,
    13.34904515479656,
    11.10340493431545,
    11.92007380866352,
    13.23724683868024,
    13.67155087160168,
    13.10482456713882,
    11.22527658932455,
    10.82300621784216,
    10.55655299881593,
    11.93905418264083,
    10.61118378234043,
    9.381881871243422,
    12.4443902414383,
    10.79018172918848,
    12.2679289839925,
    9.873804600184599,
    11.33949641835103,
    12.26926047755127,
    13.33784879145171,
    10.33215490030048,
    12.671808382455,
    11.19262517769556,
    15.48787191038721,
    13.38839581231926,
    12.90626170728282,
    12.5276258280859,
    13.19894011309694,
    14.5247248952989,
    12.93778939155904,
    13.88548470606194,
    12.77950087836109,
    13.32780306795826,
    15.91066099610471,
    13.44500137144931
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
